$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.175.34"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.615.00"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.22"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.483"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.47"
$ws.Range("E10").Value = "  +3.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.838.00"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.615.42"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.175.64"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.90"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "199.29"
$ws.Range("E20").Value = "  +5.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.28"
$ws.Range("E21").Value = "  +2.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.51"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.03"
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("E24").Value = "  +2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.60"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("E31").Value = "  +3.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.16"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.109.18"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0154"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.508"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.793"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  +8.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.749.58"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.22"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  +9.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.55"
$ws.Range("E47").Value = "  +8.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.12"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0510"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  -0.32%  "
